$wb = $excel.ActiveWorkbook

# Update sheet "展览"
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 892
$ws1.Range("F3").Value = 4523

# Update sheet "全部类型"
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 892
$ws4.Range("F3").Value = 4523
